# Weekly update: insert a new price-record row for Damasco (Castle Brite,
# Primera) at row 33 of the "Vega Modelo de Temuco" sheet, shifting the
# existing rows 33-64 down to 34-65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 33 (pushes old rows 33..64 -> 34..65,
# carries formatting down, and extends the sheet dimension automatically).
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new weekly record.
$ws.Range("A33").Value = 10
$ws.Range("B33").Value = "Vega Modelo de Temuco"
$ws.Range("C33").Value = "La Araucanía"
$ws.Range("D33").Value = 44895
$ws.Range("E33").Value = 9
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100103
$ws.Range("H33").Value = "Frutos de hueso (carozo)"
$ws.Range("I33").Value = 100103003
$ws.Range("J33").Value = "Damasco"
$ws.Range("K33").Value = "Castle Brite"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 65
$ws.Range("N33").Value = 28000
$ws.Range("O33").Value = 28000
$ws.Range("P33").Value = 28000
$ws.Range("Q33").Value = "$/bandeja 18 kilos"
$ws.Range("R33").Value = "Provincia de Limarí"
$ws.Range("S33").Value = 1556
$ws.Range("T33").Value = 18
